$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new data row for "VIDROP 2800 I.U./ML ORAL DROPS 15 ML"
#    right before "VOLTAREN 75MG/3ML 3 AMP." (currently row 39).
# ---------------------------------------------------------------------
$ws.Rows(39).Insert()

# Copy formatting (styles, number formats, fonts, borders) from the row
# above (row 38) into the freshly inserted blank row 39.
$ws.Range("A38:Q38").Copy()
$ws.Range("A39:Q39").PasteSpecial(-4122)
$ws.Range("A39").Select()
$excel.CutCopyMode = 0

$ws.Rows(39).RowHeight = 25.5

$ws.Range("A39").Value = 33
$ws.Range("C39").Value = "VIDROP 2800 I.U./ML ORAL DROPS 15 ML"
$ws.Range("H39").Value = "4:0"
$ws.Range("L39").Value = "1"
$ws.Range("N39").Value = "26.00"
$ws.Range("P39").Value = "26.0000"
$ws.Range("Q39").Value = "1:0"

$ws.Range("A39:B39").Merge()
$ws.Range("C39:G39").Merge()
$ws.Range("H39:K39").Merge()
$ws.Range("L39:M39").Merge()
$ws.Range("N39:O39").Merge()

# ---------------------------------------------------------------------
# 2) Renumber the serial numbers (column A) of every row from the old
#    "VOLTAREN..." row through the old "كريم فيبكس الازرق" row (they
#    all shift down by one place because of the insert above).
# ---------------------------------------------------------------------
$ws.Range("A40").Value = 34
$ws.Range("A41").Value = 35
$ws.Range("A42").Value = 36
$ws.Range("A43").Value = 37
$ws.Range("A44").Value = 38
$ws.Range("A45").Value = 39
$ws.Range("A46").Value = 40
$ws.Range("A47").Value = 41
$ws.Range("A48").Value = 42

# ---------------------------------------------------------------------
# 3) Insert a new data row for "سيتي بيبي رقم 2" right before
#    "كريم فيبكس الازرق" (now row 48, after the first insert above).
# ---------------------------------------------------------------------
$ws.Rows(48).Insert()

$ws.Range("A38:Q38").Copy()
$ws.Range("A48:Q48").PasteSpecial(-4122)
$ws.Range("A48").Select()
$excel.CutCopyMode = 0

$ws.Rows(48).RowHeight = 24.75

$ws.Range("A48").Value = 42
$ws.Range("C48").Value = "سيتي بيبي رقم 2"
$ws.Range("H48").Value = "1:0"
$ws.Range("L48").Value = "0"
$ws.Range("N48").Value = "180.00"
$ws.Range("P48").Value = "180.0000"
$ws.Range("Q48").Value = "1:0"

$ws.Range("A48:B48").Merge()
$ws.Range("C48:G48").Merge()
$ws.Range("H48:K48").Merge()
$ws.Range("L48:M48").Merge()
$ws.Range("N48:O48").Merge()

# "كريم فيبكس الازرق" is now row 49 - renumber its serial number.
$ws.Range("A49").Value = 43

# ---------------------------------------------------------------------
# 4) Update the grand-total cell (now on row 50) to reflect the two
#    newly added rows (26.00 + 180.00 = 206.00 more than before).
# ---------------------------------------------------------------------
$ws.Range("P50").Value = 2278.99

# ---------------------------------------------------------------------
# 5) Bump the generated timestamp shown in the footer (now row 51).
# ---------------------------------------------------------------------
$ws.Range("A51").Value = "Saturday, 2 August, 2025 4:52 PM"
